{"js": "// Replace the 25 \"XXX\u00f7Y=\" division prompts in the practice-sheet table with\n// a new set of three-digit/one-digit division problems, per the commit's\n// regenerated numbers. Each old value is unique in the document, so an\n// exact, case-sensitive search-and-replace on each pair is unambiguous.\nconst replacements = [\n  [\"578\u00f78=\", \"959\u00f73=\"],\n  [\"501\u00f73=\", \"404\u00f78=\"],\n  [\"276\u00f77=\", \"413\u00f77=\"],\n  [\"454\u00f77=\", \"189\u00f76=\"],\n  [\"221\u00f75=\", \"700\u00f74=\"],\n  [\"778\u00f74=\", \"370\u00f77=\"],\n  [\"100\u00f73=\", \"942\u00f77=\"],\n  [\"235\u00f77=\", \"506\u00f73=\"],\n  [\"397\u00f78=\", \"221\u00f78=\"],\n  [\"123\u00f79=\", \"974\u00f75=\"],\n  [\"892\u00f77=\", \"357\u00f76=\"],\n  [\"893\u00f78=\", \"965\u00f75=\"],\n  [\"120\u00f74=\", \"642\u00f75=\"],\n  [\"185\u00f73=\", \"994\u00f76=\"],\n  [\"685\u00f78=\", \"174\u00f73=\"],\n  [\"491\u00f72=\", \"452\u00f73=\"],\n  [\"589\u00f76=\", \"631\u00f72=\"],\n  [\"240\u00f74=\", \"670\u00f77=\"],\n  [\"773\u00f73=\", \"721\u00f78=\"],\n  [\"602\u00f75=\", \"881\u00f74=\"],\n  [\"733\u00f73=\", \"533\u00f77=\"],\n  [\"870\u00f78=\", \"881\u00f76=\"],\n  [\"428\u00f76=\", \"877\u00f73=\"],\n  [\"202\u00f78=\", \"209\u00f78=\"],\n  [\"745\u00f72=\", \"185\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"XXX\u00f7Y=\" division prompts in the practice-sheet table with\n# a new set of three-digit/one-digit division problems, per the commit's\n# regenerated numbers. Each old value is unique in the document, so a\n# case-sensitive Find/Replace-all on each pair is unambiguous and only ever\n# touches the one matching cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"578\u00f78=\", \"959\u00f73=\"),\n  @(\"501\u00f73=\", \"404\u00f78=\"),\n  @(\"276\u00f77=\", \"413\u00f77=\"),\n  @(\"454\u00f77=\", \"189\u00f76=\"),\n  @(\"221\u00f75=\", \"700\u00f74=\"),\n  @(\"778\u00f74=\", \"370\u00f77=\"),\n  @(\"100\u00f73=\", \"942\u00f77=\"),\n  @(\"235\u00f77=\", \"506\u00f73=\"),\n  @(\"397\u00f78=\", \"221\u00f78=\"),\n  @(\"123\u00f79=\", \"974\u00f75=\"),\n  @(\"892\u00f77=\", \"357\u00f76=\"),\n  @(\"893\u00f78=\", \"965\u00f75=\"),\n  @(\"120\u00f74=\", \"642\u00f75=\"),\n  @(\"185\u00f73=\", \"994\u00f76=\"),\n  @(\"685\u00f78=\", \"174\u00f73=\"),\n  @(\"491\u00f72=\", \"452\u00f73=\"),\n  @(\"589\u00f76=\", \"631\u00f72=\"),\n  @(\"240\u00f74=\", \"670\u00f77=\"),\n  @(\"773\u00f73=\", \"721\u00f78=\"),\n  @(\"602\u00f75=\", \"881\u00f74=\"),\n  @(\"733\u00f73=\", \"533\u00f77=\"),\n  @(\"870\u00f78=\", \"881\u00f76=\"),\n  @(\"428\u00f76=\", \"877\u00f73=\"),\n  @(\"202\u00f78=\", \"209\u00f78=\"),\n  @(\"745\u00f72=\", \"185\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
